$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 6.9
$ws.Range("C4").Value = 25
$ws.Range("B5").Value = 0.95

$ws.Range("C8").Select()
